# Highlight a set of survey-item paragraphs in yellow.
# wdYellow = 7 (WdColorIndex); setting highlight via .Font.HighlightColorIndex
# (rather than .HighlightColorIndex directly) applies it both to the runs in
# the paragraph AND to the paragraph-mark run properties (w:pPr/w:rPr), which
# mirrors how Word applies highlighting when a whole paragraph (incl. its
# end-of-paragraph mark) is selected.

$wdYellow = 7

$targets = @(
    "general public",
    "Police officers should have access to naloxone/NARCAN.",
    "People who use drugs should have access to safe inhalation supplies (glass stems and pipes).",
    "Medications used to treat addiction (buprenorphine, naltrexone, or methadone) are an appropriate treatment option for people who use drugs.",
    "Sobriety should not be a requirement to access public housing.",
    "It should be legal for adults to purchase drugs from a dispensary/shop.",
    "People who use drugs should be treated with respect.",
    "Racism",
    "Gender-based discrimination",
    "Some ways of using drugs are safer than others.",
    "People in recovery from drug use",
    "Relapse may be a part of the recovery process.",
    "People who use drugs should be forced into treatment.",
    "Using drugs is immoral.",
    "Harm reduction complements traditional addiction prevention, treatment, and recovery services.",
    "People will use more drugs if it is safer.",
    "Drug use will always be part of society.",
    "Chaotic drug use is a rational response to experiences like trauma, homelessness, hunger, and poverty.",
    "People who use drugs should be able to use medications used to treat addiction (buprenorphine, naltrexone, or methadone) for any length of time."
)

$d = $word.ActiveDocument

foreach ($target in $targets) {
    $r = $d.Content
    $found = $r.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $para = $r.Paragraphs(1)
        $para.Range.Font.HighlightColorIndex = $wdYellow
    }
}
